# Edit script: apply NP_QTR_FIN.xlsx quarterly-update diff.
# The workbook has a single sheet "NP" with sections: Income Statement, Balance
# Sheet, and Cash Flow Statement, each with a row of period-ending dates in
# columns D:K (10 quarters) as of the "before" state. The update adds two new,
# more-recent quarters as the new leftmost data columns (D, E), pushing all
# existing quarters two columns to the right (old D->F ... old K->M).
#
# Strategy:
#  1. Insert two blank columns at D (EntireColumn.Insert shifts D:K -> F:M).
#  2. Copy the number formatting from column F (which now holds what used to
#     be column D) into the two new D:E columns, row by row range, so the
#     newly inserted cells pick up the same date / number style as their
#     neighbours (Excel's native column-insert carries the left column's
#     format, but we paste explicitly to be robust).
#  3. Write the new quarter's values into D and E for every row that carries
#     data.
#  4. Apply the one genuine data correction that is not just a shift: H72
#     (Retained Earnings, one of the older quarters) changes from "NA" to
#     235700.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column D (shifts existing D:K to F:M).
$ws.Range("D1:E1").EntireColumn.Insert()

# 2. Carry over number formatting (date / #,##0 styles) from column F (the
#    former column D) into the freshly inserted D:E columns.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the two new quarters of data.

$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 240900
$ws.Range("E8").Value2 = 256200
$ws.Range("D9").Value2 = 206300
$ws.Range("E9").Value2 = 214900
$ws.Range("D10").Value2 = 34600
$ws.Range("E10").Value2 = 41300
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = -4100
$ws.Range("E14").Value2 = 700
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 223100
$ws.Range("E17").Value2 = 239700
$ws.Range("D18").Value2 = 17800
$ws.Range("E18").Value2 = 16500
$ws.Range("D20").Value2 = -3200
$ws.Range("E20").Value2 = -3200
$ws.Range("D21").Value2 = 23400
$ws.Range("E21").Value2 = 22100
$ws.Range("D22").Value2 = "NA"
$ws.Range("E22").Value2 = "NA"
$ws.Range("D23").Value2 = 14600
$ws.Range("E23").Value2 = 13300
$ws.Range("D24").Value2 = 2100
$ws.Range("E24").Value2 = 400
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 12500
$ws.Range("E26").Value2 = 12900
$ws.Range("D27").Value2 = 12600
$ws.Range("E27").Value2 = 12800
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 400
$ws.Range("E29").Value2 = -800
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 3200
$ws.Range("E32").Value2 = 3200
$ws.Range("D33").Value2 = 13000
$ws.Range("E33").Value2 = 12000
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 13000
$ws.Range("E35").Value2 = 12000
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 9900
$ws.Range("E41").Value2 = 7400
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 114800
$ws.Range("E43").Value2 = 130400
$ws.Range("D44").Value2 = 131600
$ws.Range("E44").Value2 = 137800
$ws.Range("D45").Value2 = 21600
$ws.Range("E45").Value2 = 20900
$ws.Range("D46").Value2 = 277900
$ws.Range("E46").Value2 = 296500
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 396200
$ws.Range("E48").Value2 = 397800
$ws.Range("D49").Value2 = 154700
$ws.Range("E49").Value2 = 157200
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 32400
$ws.Range("E52").Value2 = 33000
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 861200
$ws.Range("E54").Value2 = 884500
$ws.Range("D57").Value2 = 63300
$ws.Range("E57").Value2 = 69800
$ws.Range("D58").Value2 = 2300
$ws.Range("E58").Value2 = 2000
$ws.Range("D59").Value2 = 55200
$ws.Range("E59").Value2 = 60500
$ws.Range("D60").Value2 = 120800
$ws.Range("E60").Value2 = 132300
$ws.Range("D61").Value2 = 236800
$ws.Range("E61").Value2 = 247600
$ws.Range("D62").Value2 = 113400
$ws.Range("E62").Value2 = 105500
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 471000
$ws.Range("E66").Value2 = 485400
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 243200
$ws.Range("E72").Value2 = "NA"
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 390200
$ws.Range("E76").Value2 = 399100
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 13000
$ws.Range("E81").Value2 = 12000
$ws.Range("D83").Value2 = 8800
$ws.Range("E83").Value2 = 8800
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 29000
$ws.Range("E89").Value2 = 23900
$ws.Range("D91").Value2 = -10000
$ws.Range("E91").Value2 = -12300
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -5400
$ws.Range("E94").Value2 = -12900
$ws.Range("D96").Value2 = -7000
$ws.Range("E96").Value2 = -6900
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -21000
$ws.Range("E100").Value2 = -10700
$ws.Range("D101").Value2 = -100
$ws.Range("E101").Value2 = -100
$ws.Range("D102").Value2 = 2500
$ws.Range("E102").Value2 = 200

# 4. One-off data correction on a shifted cell (Retained Earnings, the
#    quarter that is now in column H): was "NA", corrected to 235700.
$ws.Range("H72").Value2 = 235700

